$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2 through 427) from 45190 to 45192.
for ($r = 2; $r -le 427; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# Row 427 picks up an explicit (customized) row height in the new file even
# though the value (15) matches the sheet default - force it to serialize.
$ws.Rows.Item(427).RowHeight = 15

# Append the new record as row 428.
$ws.Cells.Item(428, 1).Value = "A 44890-2023"

$ws.Cells.Item(428, 2).Value = 45190
$ws.Cells.Item(428, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(428, 3).Value = 45192
$ws.Cells.Item(428, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(428, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(428, 5).Value = "MORA"

$ws.Cells.Item(428, 7).Value = 3.4
$ws.Cells.Item(428, 8).Value = 0
$ws.Cells.Item(428, 9).Value = 0
$ws.Cells.Item(428, 10).Value = 0
$ws.Cells.Item(428, 11).Value = 0
$ws.Cells.Item(428, 12).Value = 0
$ws.Cells.Item(428, 13).Value = 0
$ws.Cells.Item(428, 14).Value = 0
$ws.Cells.Item(428, 15).Value = 0
$ws.Cells.Item(428, 16).Value = 0
$ws.Cells.Item(428, 17).Value = 0

# R428 stays an empty (no species match) cell but keeps the wrap-text style
# used throughout column R.
$ws.Cells.Item(428, 18).WrapText = $true

# Match the explicit row height used by every other data row.
$ws.Rows.Item(428).RowHeight = 15
